$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time entries for row 3 (2020-12-17 / serial 44103)
$ws.Range("B3").Value = 0.65844907407407405
$ws.Range("C3").Value = 0.86927083333333333
$ws.Range("D3").Value = 0.020833333333333332

# Add the new "Activity" note for row 3 (new shared string), matching F2's style
$ws.Range("F3").Value = "* Added nodes`n* Added node serialization`n* Added Undo/Redo (hopefully) and file change tracking`n* Added search tree for nodes`n* Researched the internals of ShaderGraph to learn how a bunch of things are done there, then reverse-engineered some of them"
$ws.Range("F3").NumberFormat = $ws.Range("F2").NumberFormat
$ws.Range("F3").WrapText = $true

# Keep the row height consistent with the rest of the sheet
$ws.Rows(3).RowHeight = 15

# Move the active selection to H9, matching the saved view state
$ws.Range("H9").Select()
